# Update the stock-screener table on Sheet1 (columns B..F, rows 2..27)
# to the new ticker lists, add a new row 27 (index 25) and expand the
# used range accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New contents for columns B,C,D,E,F across rows 2..27 ($null = blank cell)
$data = @(
  @("NSE:AGRITECH", "NSE:AEROFLEX", "NSE:COROMANDEL", "NSE:ABBOTINDIA", "NSE:COROMANDEL"),
  @("NSE:ANANTRAJ", "NSE:AETHER", $null, "NSE:ADANIENT", "NSE:IEX"),
  @("NSE:ARIES", "NSE:AHL", $null, "NSE:DRREDDY", "NSE:M&M"),
  @("NSE:ASMS", "NSE:ALEMBICLTD", $null, "NSE:HINDPETRO", "NSE:MOTHERSON"),
  @("NSE:AXISBANK", "NSE:ALLSEC", $null, "NSE:INDIGO", $null),
  @("NSE:CHAMBLFERT", "NSE:AMBER", $null, "NSE:IOC", $null),
  @("NSE:EIFFL", "NSE:APCOTEXIND", $null, "NSE:MARUTI", $null),
  @("NSE:FMGOETZE", "NSE:APOLLOTYRE", $null, "NSE:NAVINFLUOR", $null),
  @("NSE:GENUSPAPER", "NSE:ASIANPAINT", $null, $null, $null),
  @("NSE:GSFC", "NSE:CANFINHOME", $null, $null, $null),
  @("NSE:ICRA", "NSE:CEATLTD", $null, $null, $null),
  @("NSE:IEX", "NSE:ESSENTIA", $null, $null, $null),
  @("NSE:M&M", "NSE:FAZE3Q", $null, $null, $null),
  @("NSE:MADRASFERT", "NSE:GHCL", $null, $null, $null),
  @("NSE:MOTHERSON", "NSE:GICHSGFIN", $null, $null, $null),
  @("NSE:RELINFRA", "NSE:GINNIFILA", $null, $null, $null),
  @($null, "NSE:HINDPETRO", $null, $null, $null),
  @($null, "NSE:JKTYRE", $null, $null, $null),
  @($null, "NSE:KDDL", $null, $null, $null),
  @($null, "NSE:LICHSGFIN", $null, $null, $null),
  @($null, "NSE:MEDPLUS", $null, $null, $null),
  @($null, "NSE:MOVALUE", $null, $null, $null),
  @($null, "NSE:MPSLTD", $null, $null, $null),
  @($null, "NSE:NECLIFE", $null, $null, $null),
  @($null, "NSE:NUCLEUS", $null, $null, $null),
  @($null, "NSE:PFIZER", $null, $null, $null)
)

# Row 27 is brand new: give its A cell the same style (border/bold/
# alignment) as the existing numbered rows before writing any values.
$ws.Range("A26").Copy()
$ws.Range("A27").PasteSpecial(-4122)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowNum = 2 + $r
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $colNum = 2 + $c
        $val = $row[$c]
        $cell = $ws.Cells.Item($rowNum, $colNum)
        if ($val -ne $null) {
            $cell.Value = $val
        } else {
            $cell.ClearContents()
        }
    }
}

# Row 27's running index (column A) continues the 0-based sequence -> 25
$ws.Cells.Item(27, 1).Value = 25
